# Garonzi_Imran_Klein_CocomoMdS_WBS.xlsx -- "Add files via upload"
#
# This localizes the WBS header row into Italian and pushes out the "End"
# (minimum start / maximum start) dates for a few tasks whose schedule
# slipped, then widens the two date columns so the longer Italian header
# text fits comfortably.
#
#  1. Row 5 header translations:
#       B5  Level              -> Livello
#       D5  Task Description   -> Descrizione Task
#       E5  Assigned To        -> Assegnato a
#       F5  Start              -> Inizio minimo
#       G5  End                -> Inizio massimo
#  2. Column G ("End") date slips:
#       Fase 1 (rows 7:9)   14 Dec 2022 -> 15 Dec 2022
#       Fase 2 (rows 11:17) 17 Dec 2022 -> 19 Dec 2022
#       Fase 3 (rows 19:23) 18 Dec 2022 -> 22 Dec 2022
#  3. Widen columns F:G (Start/End) to fit the new header labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header translations (row 5) ---------------------------------------
$ws.Range("B5").Value2 = "Livello"
$ws.Range("D5").Value2 = "Descrizione Task"
$ws.Range("E5").Value2 = "Assegnato a"
$ws.Range("F5").Value2 = "Inizio minimo"
$ws.Range("G5").Value2 = "Inizio massimo"

# --- 2. Updated "End" dates in column G ------------------------------------
$ws.Range("G7:G9").Value2 = 44910      # Fase 1 -> 15/12/2022
$ws.Range("G11:G17").Value2 = 44914    # Fase 2 -> 19/12/2022
$ws.Range("G19:G23").Value2 = 44917    # Fase 3 -> 22/12/2022

# --- 3. Widen the Start/End columns (F:G) ----------------------------------
$ws.Range("F1:G1").EntireColumn.ColumnWidth = 13
